$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2150
$ws.Range("I4").Value = 1900
$ws.Range("J4").Value = 2400
$ws.Range("K4").Value = 1900
$ws.Range("L4").Value = 2400
$ws.Range("M4").Value = -1786
$ws.Range("N4").Value = -2628
$ws.Range("H53").Value = 5044.5
$ws.Range("I53").Value = 5044.5
$ws.Range("K53").Value = 5044.5
$ws.Range("M53").Value = -4407.5
$ws.Range("H64").Value = 2801.6667
$ws.Range("I64").Value = 2801.6667
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 2801.6667
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -2553.6667
$ws.Range("H67").Value = 2801.6667
$ws.Range("I67").Value = 2801.6667
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 2801.6667
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -1943.6667
$ws.Range("H74").Value = 3499
$ws.Range("I74").Value = 2373.5
$ws.Range("K74").Value = 2373.5
$ws.Range("M74").Value = -1437.5
$ws.Range("H77").Value = 3499
$ws.Range("I77").Value = 2373.5
$ws.Range("K77").Value = 11867.5
$ws.Range("M77").Value = -7187.5
$ws.Range("H112").Value = 1952.3334
$ws.Range("J112").Value = 1952.3334
$ws.Range("L112").Value = 5857.0002
$ws.Range("N112").Value = -8073.0002
$ws.Range("H132").Value = 1089.2307
$ws.Range("I132").Value = 970.9
$ws.Range("J132").Value = 1483.6666
$ws.Range("K132").Value = 2912.7
$ws.Range("L132").Value = 4450.9998
$ws.Range("M132").Value = -382.6999999999998
$ws.Range("N132").Value = -9510.9998
$ws.Range("H138").Value = 1607.0862
$ws.Range("J138").Value = 2248.5806
$ws.Range("L138").Value = 6745.7418
$ws.Range("N138").Value = -17025.7418

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 999
$ws.Range("I21").Value = 999
$ws.Range("K21").Value = 999
$ws.Range("M21").Value = -625
$ws.Range("H32").Value = 4970.8037
$ws.Range("I32").Value = 3609.2173
$ws.Range("K32").Value = 3609.2173
$ws.Range("M32").Value = -3322.2173
$ws.Range("H45").Value = 7501507.5
$ws.Range("I45").Value = 22500724
$ws.Range("K45").Value = 22500724
$ws.Range("M45").Value = -22500347
$ws.Range("H61").Value = 7635.5557
$ws.Range("I61").Value = 9013
$ws.Range("J61").Value = 5471
$ws.Range("K61").Value = 9013
$ws.Range("L61").Value = 5471
$ws.Range("M61").Value = -8801
$ws.Range("N61").Value = -5895
$ws.Range("H97").Value = 967.04346
$ws.Range("I97").Value = 938.3182
$ws.Range("K97").Value = 938.3182
$ws.Range("M97").Value = -442.3182
$ws.Range("H102").Value = 1056.8462
$ws.Range("I102").Value = 769.25
$ws.Range("K102").Value = 769.25
$ws.Range("M102").Value = 852.75
$ws.Range("H110").Value = 240.63637
$ws.Range("I110").Value = 240.63637
$ws.Range("K110").Value = 240.63637
$ws.Range("M110").Value = 1804.36363
$ws.Range("H122").Value = 1029.1666
$ws.Range("I122").Value = 843.13336
$ws.Range("J122").Value = 1339.2222
$ws.Range("K122").Value = 2529.40008
$ws.Range("L122").Value = 4017.6666
$ws.Range("M122").Value = -79.40008000000034
$ws.Range("N122").Value = -8917.6666
$ws.Range("H132").Value = 1557.1471
$ws.Range("I132").Value = 1201.1852
$ws.Range("K132").Value = 3603.5556
$ws.Range("M132").Value = -1073.5556
$ws.Range("H136").Value = 7635.5557
$ws.Range("I136").Value = 9013
$ws.Range("J136").Value = 5471
$ws.Range("K136").Value = 27039
$ws.Range("L136").Value = 16413
$ws.Range("M136").Value = -24489
$ws.Range("N136").Value = -21513

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1997
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H105").Value = 2226.1785
$ws.Range("J105").Value = 3692.5
$ws.Range("L105").Value = 3692.5
$ws.Range("N105").Value = -7186.5
$ws.Range("H107").Value = 995.8461
$ws.Range("I107").Value = 655.4
$ws.Range("K107").Value = 655.4
$ws.Range("M107").Value = 1264.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2272.5908
$ws.Range("I31").Value = 2214.1667
$ws.Range("K31").Value = 2214.1667
$ws.Range("M31").Value = -1919.1667
$ws.Range("H34").Value = 2272.5908
$ws.Range("I34").Value = 2214.1667
$ws.Range("K34").Value = 2214.1667
$ws.Range("M34").Value = -2012.1667
$ws.Range("H43").Value = 30000
$ws.Range("J43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30368
$ws.Range("H99").Value = 1001907.7
$ws.Range("I99").Value = 1251822.1
$ws.Range("J99").Value = 2250
$ws.Range("K99").Value = 1251822.1
$ws.Range("L99").Value = 2250
$ws.Range("M99").Value = -1250324.1
$ws.Range("N99").Value = -5246
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H126").Value = 1001907.7
$ws.Range("I126").Value = 1251822.1
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 3755466.3
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = -3752996.3
$ws.Range("N126").Value = -11690
$ws.Range("H132").Value = 2782.889
$ws.Range("J132").Value = 4177.25
$ws.Range("L132").Value = 12531.75
$ws.Range("N132").Value = -17591.75
$ws.Range("H134").Value = 2659
$ws.Range("I134").Value = 2269.0667
$ws.Range("K134").Value = 6807.2001
$ws.Range("M134").Value = -4272.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2165.2354
$ws.Range("I114").Value = 325.4
$ws.Range("J114").Value = 2931.8333
$ws.Range("K114").Value = 976.1999999999999
$ws.Range("L114").Value = 8795.499899999999
$ws.Range("M114").Value = 2277.8
$ws.Range("N114").Value = -15303.4999
$ws.Range("H129").Value = 104269.14
$ws.Range("J129").Value = 181789
$ws.Range("L129").Value = 545367
$ws.Range("N129").Value = -555367
$ws.Range("H131").Value = 15318.75
$ws.Range("I131").Value = 840
$ws.Range("J131").Value = 15948.261
$ws.Range("K131").Value = 2520
$ws.Range("L131").Value = 47844.783
$ws.Range("M131").Value = 2520
$ws.Range("N131").Value = -57924.783
$ws.Range("H139").Value = 5048.423
$ws.Range("I139").Value = 6143.579
$ws.Range("J139").Value = 2075.8572
$ws.Range("K139").Value = 18430.737
$ws.Range("L139").Value = 6227.571599999999
$ws.Range("M139").Value = -13290.737
$ws.Range("N139").Value = -16507.5716
$ws.Range("H140").Value = 2725
$ws.Range("I140").Value = 1189.75
$ws.Range("K140").Value = 3569.25
$ws.Range("M140").Value = 1610.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 108
$ws.Range("I2").Value = 83
$ws.Range("J2").Value = 133
$ws.Range("K2").Value = 83
$ws.Range("L2").Value = 133
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = -359
$ws.Range("H102").Value = 2677.6667
$ws.Range("I102").Value = 2523.3845
$ws.Range("J102").Value = 3078.8
$ws.Range("K102").Value = 2523.3845
$ws.Range("L102").Value = 3078.8
$ws.Range("M102").Value = -901.3845000000001
$ws.Range("N102").Value = -6322.8
$ws.Range("H113").Value = 1336.4
$ws.Range("I113").Value = 1085.5
$ws.Range("K113").Value = 1085.5
$ws.Range("M113").Value = 1084.5
$ws.Range("H122").Value = 1415.7727
$ws.Range("I122").Value = 1455.2142
$ws.Range("J122").Value = 1346.75
$ws.Range("K122").Value = 4365.642599999999
$ws.Range("L122").Value = 4040.25
$ws.Range("M122").Value = -1915.642599999999
$ws.Range("N122").Value = -8940.25
$ws.Range("H132").Value = 3208924.8
$ws.Range("I132").Value = 5496871.5
$ws.Range("K132").Value = 16490614.5
$ws.Range("M132").Value = -16488084.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 529.05554
$ws.Range("I55").Value = 502.4
$ws.Range("K55").Value = 502.4
$ws.Range("M55").Value = -329.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1355.1428
$ws.Range("J96").Value = 1536.6
$ws.Range("L96").Value = 1536.6
$ws.Range("N96").Value = -4282.6
$ws.Range("H113").Value = 505.32144
$ws.Range("J113").Value = 654.46155
$ws.Range("L113").Value = 1963.38465
$ws.Range("N113").Value = -6303.38465
$ws.Range("H132").Value = 1143.6531
$ws.Range("I132").Value = 797.7954999999999
$ws.Range("K132").Value = 2393.3865
$ws.Range("M132").Value = 136.6135000000004
